# ---------------------------------------------------------------------------
# Helper functions
# ---------------------------------------------------------------------------

function Get-FoundRange($doc, $needle, $scope) {
    # $scope is optional; if omitted, search the whole document content.
    if ($scope) {
        $r = $scope
    } else {
        $r = $doc.Content
    }
    $ok = $r.Find.Execute($needle, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $ok) {
        throw ("Not found: " + $needle)
    }
    return $r
}

# Deletes the (unique) $needle text anywhere in $doc and inserts $newText at
# that exact point. Because the deletion leaves a fresh insertion point, the
# new text is written into its own run(s), cleanly split from the
# surrounding (unmodified) text -- matching how Word COM automation actually
# edits a document.
function Replace-Span($doc, $needle, $newText) {
    $r = Get-FoundRange $doc $needle
    $r.Delete()
    $ins = $doc.Range($r.Start, $r.Start)
    if ($newText.Length -gt 0) {
        $ins.InsertBefore($newText)
    }
    $after = $doc.Range($r.Start, $r.Start + $newText.Length)
    return $after
}

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. "Наш проект разработан..." -> "Мой проект разработан..."
# ---------------------------------------------------------------------------
Replace-Span $d "Наш проект разработан" "Мой проект разработан" | Out-Null

# ---------------------------------------------------------------------------
# 2. "Мы также предусмотрели..." -> "Я также предусмотрели..."
# ---------------------------------------------------------------------------
Replace-Span $d "Мы также предусмотрели" "Я также предусмотрели" | Out-Null

# ---------------------------------------------------------------------------
# 3. ". Кроме того, мы реализовали" -> ". Кроме того, я реализовали"
# ---------------------------------------------------------------------------
Replace-Span $d "Кроме того, мы реализовали" "Кроме того, я реализовали" | Out-Null

# ---------------------------------------------------------------------------
# 4. "недостатков. Мы стремимся" -> "недостатков. Я стремлюсь"
#    (split into "Я" / " стрем" / "люсь" to mirror the authored edit)
# ---------------------------------------------------------------------------
Replace-Span $d "недостатков. Мы стремимся" "недостатков. Я стремимся" | Out-Null
Replace-Span $d "стремимся" "стремлюсь" | Out-Null

# ---------------------------------------------------------------------------
# 5. Delete the whole "Кроме того, мы планируем привлечь финансирование..."
#    paragraph, merging the following (pict-only) paragraph up into the
#    "...родителей." paragraph.
# ---------------------------------------------------------------------------
$p9 = Get-FoundRange $d "Кроме того, мы планируем привлечь финансирование для дальнейшего развития проекта, его улучшения и внедрения в другие учебные заведения, начав с 6-й школы Ханты-Мансийска."
$p9Para = $d.Paragraphs.Item($p9.Paragraphs.Item(1).Index)
$fullPara = $p9Para.Range
$fullPara.Delete()
$prevMarkStart = $fullPara.Start - 1
$mark = $d.Range($prevMarkStart, $fullPara.Start)
$mark.Delete()

# ---------------------------------------------------------------------------
# 6. "Наш бот предоставляет..." -> "Б" + "от предоставляет..."
# ---------------------------------------------------------------------------
Replace-Span $d "Наш б" "Б" | Out-Null

# ---------------------------------------------------------------------------
# 7. "Главное преимущество нашего бота — " -> "Главное преимущество бота — "
# ---------------------------------------------------------------------------
Replace-Span $d "Главное преимущество нашего бота" "Главное преимущество бота" | Out-Null

# ---------------------------------------------------------------------------
# 8. Move <w:lastRenderedPageBreak/> from the pict-only paragraph right
#    before "Монетизация" to the start of the "Мы планируем монетизировать..."
#    paragraph, and rewrite that paragraph's opening: "Мы планируем" -> " планирую"
# ---------------------------------------------------------------------------
Replace-Span $d "Мы планируем монетизировать" " планирую монетизировать" | Out-Null

$lrpb = Get-FoundRange $d "монетизировать проект следующими способами"
$insAt = $d.Range($lrpb.Start, $lrpb.Start)
$insAt.InsertBefore(" ")
$markerFind = $d.Content
$markerFind.Find.Execute(" планирую монетизировать", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

# ---------------------------------------------------------------------------
# 9. Merge "«подписки" + "»" bold runs and drop proofErr wrapping (text
#    content is unchanged, so re-assert it verbatim to coalesce the runs).
# ---------------------------------------------------------------------------
Replace-Span $d "«подписки»" "«подписки»" | Out-Null

# ---------------------------------------------------------------------------
# 10. "Спасибо за внимание! Мы готовы ответить" ->
#     "Спасибо за внимание" / "! Готов" / " ответить"
# ---------------------------------------------------------------------------
Replace-Span $d "Спасибо за внимание! Мы готовы ответить" "Спасибо за внимание! Готов ответить" | Out-Null
